$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(3,9,10,13,18,20,23,27,29,33,39,41,46,53,54,60,61,63,70,72,73,74,76,78,84,86,90,96,97,98,101,107,108,110,112,118,119,120,124,130,133,134,137,140,144,151,154,159,162,164,165,166,171,172,173,176,177,178,182,183,185,186,187,192)
$cols = @("G","H","I","J","K","L","M","N")

foreach ($r in $rows) {
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = "N/A"
    }
}
